$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.487945148329342
$ws.Range("B2").Value = -4.710020033330782

$ws.Range("A3").Value = -0.4841097597569224
$ws.Range("B3").Value = 0.5586593429552213

$ws.Range("A4").Value = 0.8445102526334978
$ws.Range("B4").Value = -2.79945871415623

$ws.Range("A5").Value = 0.7452019343263052
$ws.Range("B5").Value = 0.6458486690017327

$ws.Range("A6").Value = -0.8194715288099588
$ws.Range("B6").Value = -2.363486361546164

$ws.Range("A7").Value = -0.08620687001602649
$ws.Range("B7").Value = -0.6280732378937954

$ws.Range("A8").Value = 0.7993982338251456
$ws.Range("B8").Value = 0.8574661753142296
